$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.503233333333334
$ws.Cells.Item(2, 8).Value = 7.5097
$ws.Cells.Item(2, 9).Value = 0.9574104874676208
$ws.Cells.Item(2, 10).Value = 0.9574104874676207
$ws.Cells.Item(2, 13).Value = 449.3583473333333
$ws.Cells.Item(2, 14).Value = 1348.075042
$ws.Cells.Item(2, 15).Value = 0.959704436884883
$ws.Cells.Item(2, 16).Value = 0.9597044368848828
$ws.Cells.Item(2, 17).Value = 1124.848793656378
$ws.Cells.Item(2, 18).Value = 10123.6391429074
$ws.Cells.Item(2, 19).Value = 0.9188310927427944
$ws.Cells.Item(2, 20).Value = 0.9188310927427941
$ws.Cells.Item(3, 7).Value = 2.503233333333334
$ws.Cells.Item(3, 8).Value = 7.5097
$ws.Cells.Item(3, 9).Value = 0.9574104874676208
$ws.Cells.Item(3, 10).Value = 0.9574104874676207
$ws.Cells.Item(3, 15).Value = 0.01202662913387072
$ws.Cells.Item(3, 16).Value = 0.01202662913387072
$ws.Cells.Item(3, 17).Value = 14.09615164112222
$ws.Cells.Item(3, 18).Value = 126.8653647701
$ws.Cells.Item(3, 19).Value = 0.01151442086165145
$ws.Cells.Item(3, 20).Value = 0.01151442086165145
$ws.Cells.Item(4, 7).Value = 2.503233333333334
$ws.Cells.Item(4, 8).Value = 7.5097
$ws.Cells.Item(4, 9).Value = 0.9574104874676208
$ws.Cells.Item(4, 10).Value = 0.9574104874676207
$ws.Cells.Item(4, 13).Value = 7.708291333333332
$ws.Cells.Item(4, 14).Value = 23.124874
$ws.Cells.Item(4, 15).Value = 0.01646276615823874
$ws.Cells.Item(4, 16).Value = 0.01646276615823874
$ws.Cells.Item(4, 17).Value = 19.29565180864444
$ws.Cells.Item(4, 18).Value = 173.6608662778
$ws.Cells.Item(4, 19).Value = 0.01576162497262481
$ws.Cells.Item(4, 20).Value = 0.0157616249726248
$ws.Cells.Item(5, 7).Value = 2.503233333333334
$ws.Cells.Item(5, 8).Value = 7.5097
$ws.Cells.Item(5, 9).Value = 0.9574104874676208
$ws.Cells.Item(5, 10).Value = 0.9574104874676207
$ws.Cells.Item(5, 13).Value = 1.356257333333333
$ws.Cells.Item(5, 14).Value = 4.068772
$ws.Cells.Item(5, 15).Value = 0.002896588408965574
$ws.Cells.Item(5, 16).Value = 0.002896588408965573
$ws.Cells.Item(5, 17).Value = 3.395028565377778
$ws.Cells.Item(5, 18).Value = 30.5552570884
$ws.Cells.Item(5, 19).Value = 0.002773224120620791
$ws.Cells.Item(5, 20).Value = 0.00277322412062079
$ws.Cells.Item(6, 7).Value = 2.503233333333334
$ws.Cells.Item(6, 8).Value = 7.5097
$ws.Cells.Item(6, 9).Value = 0.9574104874676208
$ws.Cells.Item(6, 10).Value = 0.9574104874676207
$ws.Cells.Item(6, 13).Value = 4.171694666666667
$ws.Cells.Item(6, 14).Value = 12.515084
$ws.Cells.Item(6, 15).Value = 0.008909579414042005
$ws.Cells.Item(6, 16).Value = 0.008909579414042003
$ws.Cells.Item(6, 17).Value = 10.44272514608889
$ws.Cells.Item(6, 18).Value = 93.98452631480002
$ws.Cells.Item(6, 19).Value = 0.008530124769929434
$ws.Cells.Item(6, 20).Value = 0.008530124769929433
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.111354
$ws.Cells.Item(7, 8).Value = 0.334062
$ws.Cells.Item(7, 9).Value = 0.04258951253237923
$ws.Cells.Item(7, 10).Value = 0.04258951253237923
$ws.Cells.Item(7, 13).Value = 449.3583473333333
$ws.Cells.Item(7, 14).Value = 1348.075042
$ws.Cells.Item(7, 15).Value = 0.959704436884883
$ws.Cells.Item(7, 16).Value = 0.9597044368848828
$ws.Cells.Item(7, 17).Value = 50.03784940895599
$ws.Cells.Item(7, 18).Value = 450.3406446806039
$ws.Cells.Item(7, 19).Value = 0.04087334414208867
$ws.Cells.Item(7, 20).Value = 0.04087334414208867
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.111354
$ws.Cells.Item(8, 8).Value = 0.334062
$ws.Cells.Item(8, 9).Value = 0.04258951253237923
$ws.Cells.Item(8, 10).Value = 0.04258951253237923
$ws.Cells.Item(8, 15).Value = 0.01202662913387072
$ws.Cells.Item(8, 16).Value = 0.01202662913387072
$ws.Cells.Item(8, 17).Value = 0.627054157894
$ws.Cells.Item(8, 18).Value = 5.643487421046
$ws.Cells.Item(8, 19).Value = 0.0005122082722192641
$ws.Cells.Item(8, 20).Value = 0.0005122082722192639
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.111354
$ws.Cells.Item(9, 8).Value = 0.334062
$ws.Cells.Item(9, 9).Value = 0.04258951253237923
$ws.Cells.Item(9, 10).Value = 0.04258951253237923
$ws.Cells.Item(9, 13).Value = 7.708291333333332
$ws.Cells.Item(9, 14).Value = 23.124874
$ws.Cells.Item(9, 15).Value = 0.01646276615823874
$ws.Cells.Item(9, 16).Value = 0.01646276615823874
$ws.Cells.Item(9, 17).Value = 0.8583490731319998
$ws.Cells.Item(9, 18).Value = 7.725141658187999
$ws.Cells.Item(9, 19).Value = 0.0007011411856139375
$ws.Cells.Item(9, 20).Value = 0.0007011411856139374
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.111354
$ws.Cells.Item(10, 8).Value = 0.334062
$ws.Cells.Item(10, 9).Value = 0.04258951253237923
$ws.Cells.Item(10, 10).Value = 0.04258951253237923
$ws.Cells.Item(10, 13).Value = 1.356257333333333
$ws.Cells.Item(10, 14).Value = 4.068772
$ws.Cells.Item(10, 15).Value = 0.002896588408965574
$ws.Cells.Item(10, 16).Value = 0.002896588408965573
$ws.Cells.Item(10, 17).Value = 0.151024679096
$ws.Cells.Item(10, 18).Value = 1.359222111864
$ws.Cells.Item(10, 19).Value = 0.0001233642883447837
$ws.Cells.Item(10, 20).Value = 0.0001233642883447837
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.111354
$ws.Cells.Item(11, 8).Value = 0.334062
$ws.Cells.Item(11, 9).Value = 0.04258951253237923
$ws.Cells.Item(11, 10).Value = 0.04258951253237923
$ws.Cells.Item(11, 13).Value = 4.171694666666667
$ws.Cells.Item(11, 14).Value = 12.515084
$ws.Cells.Item(11, 15).Value = 0.008909579414042005
$ws.Cells.Item(11, 16).Value = 0.008909579414042003
$ws.Cells.Item(11, 17).Value = 0.464534887912
$ws.Cells.Item(11, 18).Value = 4.180813991208
$ws.Cells.Item(11, 19).Value = 0.0003794546441125699
$ws.Cells.Item(11, 20).Value = 0.0003794546441125699

$wb.Save()
